# Update the "Price" (column D) values for several cryptocurrency rows
# as per the latest data refresh (GitHub Actions symbol-list update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new price text (must stay as text, matching the
# original inline-string cell type and preserving exact decimal formatting).
$updates = [ordered]@{
    2  = "245.64"
    4  = "5.346"
    5  = "0.05837"
    6  = "6.484"
    7  = "3.365"
    8  = "0.8118"
    9  = "0.9226"
    10 = "0.1412"
    11 = "0.07399"
    12 = "0.03103"
    13 = "0.03057"
    14 = "0.09373"
    15 = "3.865"
    16 = "0.001560"
    17 = "0.04702"
    18 = "0.0006045"
    19 = "0.006083"
    20 = "0.001248"
    21 = "0.004693"
    22 = "0.00008817"
    23 = "3.593"
    28 = "0.0002659"
    40 = "0.03850"
    41 = "0.006413"
    43 = "0.002675"
    44 = "0.008572"
    45 = "0.00005262"
    47 = "0.6536"
    48 = "0.001862"
    49 = "0.00002105"
    50 = "0.0002005"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    # Force text storage so trailing/leading zeros in the decimal text are preserved
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
}
